$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "audioBGMType" column (C) and re-key rows by id values.
# New layout: only id (A) and name (B) columns remain.
# Clear column C entirely (data + header).
$ws.Range("C1:C5").Clear()

# Rewrite header row.
$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "name"

# Rewrite data rows: id becomes the numeric BGM key (1000-1003),
# name keeps the BGM_Boss_* / BGM_4 values, reordered so BGM_4 is last.
$ws.Range("A2").Value = 1000
$ws.Range("B2").Value = "BGM_Boss_01"

$ws.Range("A3").Value = 1001
$ws.Range("B3").Value = "BGM_Boss_02"

$ws.Range("A4").Value = 1002
$ws.Range("B4").Value = "BGM_Boss_03"

$ws.Range("A5").Value = 1003
$ws.Range("B5").Value = "BGM_4"

# Update the selection to match the post-edit cursor position.
$ws.Range("L25").Select()
